$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current "Tipo" column (D) so it becomes E,
# and the new column D will hold "MAE" values.
$ws.Columns.Item(4).Insert()

# Header row
$ws.Range("D1").Value = "MAE"

# Update existing MSE (B) and R2 (C) values, and set new MAE (D) values per row.
$ws.Range("B2").Value = 1.327087135759158
$ws.Range("C2").Value = 0.7838038071973591
$ws.Range("D2").Value = 0.9503225238534604

$ws.Range("B3").Value = 6.825308633341774
$ws.Range("C3").Value = 0.9024950065214221
$ws.Range("D3").Value = 1.950878830649375

$ws.Range("B4").Value = 4.1929622691235
$ws.Range("C4").Value = 0.7928846540915259
$ws.Range("D4").Value = 1.623828393580852

$ws.Range("B5").Value = 2.978825768303484
$ws.Range("C5").Value = 0.9980433839655243
$ws.Range("D5").Value = 1.270707098373362

$ws.Range("B6").Value = 2.109471422953313
$ws.Range("C6").Value = 0.9770136838027224
$ws.Range("D6").Value = 1.159455959981469

$ws.Range("B7").Value = 1.797825224332745
$ws.Range("C7").Value = 0.9989920155773655
$ws.Range("D7").Value = 1.051486986069193

$ws.Range("B8").Value = 2.27808364023073
$ws.Range("C8").Value = 0.9973905616280521
$ws.Range("D8").Value = 1.246247103210239

$ws.Range("B9").Value = 15.83496996529479
$ws.Range("C9").Value = 0.8107038761795236
$ws.Range("D9").Value = 3.307229920429935

$ws.Range("B10").Value = 1.562186803622715
$ws.Range("C10").Value = 0.9953431695246219
$ws.Range("D10").Value = 1.005006491178466
